# Insert two new price-record rows into the weekly "Fruta / hortaliza" data
# table of the "Vega Modelo de Temuco - Frutilla" sheet.
#
# The table already held rows 2..277 of data (row 1 is the header). Two new
# rows are inserted right before the former row 172, which pushes every
# following row down by two (the former row 172 becomes row 174, ..., the
# former row 277 becomes row 279).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above row 172 (existing rows 172:277 shift to 174:279).
$ws.Rows("172:173").Insert()

# --- New row 172 -----------------------------------------------------------
$ws.Cells.Item(172, 1).Value  = 10
$ws.Cells.Item(172, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(172, 3).Value  = "La Araucanía"
$ws.Cells.Item(172, 4).Value  = 44873
$ws.Cells.Item(172, 5).Value  = 9
$ws.Cells.Item(172, 6).Value  = "Fruta"
$ws.Cells.Item(172, 7).Value  = 100101
$ws.Cells.Item(172, 8).Value  = "Berries"
$ws.Cells.Item(172, 9).Value  = 100112025
$ws.Cells.Item(172, 10).Value = "Frutilla"
$ws.Cells.Item(172, 11).Value = "Sin especificar"
$ws.Cells.Item(172, 12).Value = "Primera"
$ws.Cells.Item(172, 13).Value = 1200
$ws.Cells.Item(172, 14).Value = 9000
$ws.Cells.Item(172, 15).Value = 9500
$ws.Cells.Item(172, 16).Value = 9333
$ws.Cells.Item(172, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(172, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(172, 19).Value = 1333
$ws.Cells.Item(172, 20).Value = 7

# --- New row 173 -----------------------------------------------------------
$ws.Cells.Item(173, 1).Value  = 10
$ws.Cells.Item(173, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(173, 3).Value  = "La Araucanía"
$ws.Cells.Item(173, 4).Value  = 44873
$ws.Cells.Item(173, 5).Value  = 9
$ws.Cells.Item(173, 6).Value  = "Fruta"
$ws.Cells.Item(173, 7).Value  = 100101
$ws.Cells.Item(173, 8).Value  = "Berries"
$ws.Cells.Item(173, 9).Value  = 100112025
$ws.Cells.Item(173, 10).Value = "Frutilla"
$ws.Cells.Item(173, 11).Value = "Sin especificar"
$ws.Cells.Item(173, 12).Value = "Tercera"
$ws.Cells.Item(173, 13).Value = 130
$ws.Cells.Item(173, 14).Value = 6000
$ws.Cells.Item(173, 15).Value = 6500
$ws.Cells.Item(173, 16).Value = 6308
$ws.Cells.Item(173, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(173, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(173, 19).Value = 901
$ws.Cells.Item(173, 20).Value = 7
